$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the empty paragraph right after
#    the phone number (it is not re-added anywhere else in this area).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Delete the old "Personal Website" block:
#      - empty italic paragraph
#      - "Personal Website (www.samblackmore.xyz)" paragraph
#      - "Created portfolio website ..." bullet paragraph
#    and re-create the "_GoBack" bookmark on the empty bold paragraph that is
#    left behind (the one that used to follow the deleted block), matching
#    the paragraph that originally carried the bookmark before this block of
#    text existed.
# ---------------------------------------------------------------------------
$websiteHeader = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Personal Website (www.samblackmore.xyz)") {
        $websiteHeader = $para
        break
    }
}

if ($websiteHeader -ne $null) {
    $emptyItalicPara = $websiteHeader.Previous()
    $bulletPara = $websiteHeader.Next()

    $deleteStart = $emptyItalicPara.Range.Start
    $deleteEnd = $bulletPara.Range.End

    $delRange = $d.Range($deleteStart, $deleteEnd)
    $delRange.Delete()

    # After the delete, the paragraph that used to trail the removed block
    # (the empty bold paragraph) now begins exactly at $deleteStart; re-fetch
    # it fresh (stale Paragraph/Range object references captured before a
    # mutation do not track the shift) and stamp the bookmark on it.
    $trailingEmptyRange = $d.Range($deleteStart, $deleteStart)
    $trailingEmptyPara = $trailingEmptyRange.Paragraphs.Item(1)
    $d.Bookmarks.Add("_GoBack", $trailingEmptyPara.Range)
}

# ---------------------------------------------------------------------------
# 3. Drop the stale "lastRenderedPageBreak" hint that sits in front of the
#    "Music that I've created ..." sentence (the bullet under "Music").
# ---------------------------------------------------------------------------
$musicHeader = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Music") {
        $musicHeader = $para
        break
    }
}

if ($musicHeader -ne $null) {
    $musicBullet = $musicHeader.Next()
    $bulletStart = $musicBullet.Range.Start
    $firstChar = $d.Range($bulletStart, $bulletStart + 1)
    $firstCharText = $firstChar.Text
    $firstChar.Delete()
    $reinsertPoint = $d.Range($bulletStart, $bulletStart)
    $reinsertPoint.InsertBefore($firstCharText)
}
